$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.06847097740271614
$ws.Range("C2").Value = 0.6620077945808225
$ws.Range("D2").Value = 0.8576394255630952
$ws.Range("E2").Value = 0.9260882385405266
$ws.Range("F2").Value = 0.9347485423480986
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = 0.1628294350496699
$ws.Range("C3").Value = 0.6137740522511609
$ws.Range("D3").Value = 0.7803318254486972
$ws.Range("E3").Value = 0.8833639258248535
$ws.Range("F3").Value = 0.8790130036258593
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = -0.06854388644398214
$ws.Range("C4").Value = 0.6834144119998282
$ws.Range("D4").Value = 0.8849781572300088
$ws.Range("E4").Value = 0.9407327767384364
$ws.Range("F4").Value = 0.9501848094793431
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = 0.1479423124587954
$ws.Range("C5").Value = 0.6621921511100777
$ws.Range("D5").Value = 0.852793974404251
$ws.Range("E5").Value = 0.9234684479744021
$ws.Range("F5").Value = 0.9234571032261442
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = -0.08040621195020153
$ws.Range("C6").Value = 0.6834744100144207
$ws.Range("D6").Value = 0.8476544914969599
$ws.Range("E6").Value = 0.920681536415801
$ws.Range("F6").Value = 0.9294752172076342
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.1165699234898814
$ws.Range("C7").Value = 0.6791032791661118
$ws.Range("D7").Value = 0.8276536103283203
$ws.Range("E7").Value = 0.9097546978874692
$ws.Range("F7").Value = 0.9147010340487955
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = -0.1318557986420442
$ws.Range("C8").Value = 0.6256777560484761
$ws.Range("D8").Value = 0.7095626088736116
$ws.Range("E8").Value = 0.8423553934495888
$ws.Range("F8").Value = 0.843773152665554
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.04870157088247548
$ws.Range("C9").Value = 0.6198098084369814
$ws.Range("D9").Value = 0.7254827416971069
$ws.Range("E9").Value = 0.8517527468092527
$ws.Range("F9").Value = 0.8627739369610303
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = -0.06798252598835164
$ws.Range("C10").Value = 0.5847756573626377
$ws.Range("D10").Value = 0.6476270503415246
$ws.Range("E10").Value = 0.8047527883403229
$ws.Range("F10").Value = 0.8139351567699135
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.01270493608304687
$ws.Range("C11").Value = 0.6241980251273542
$ws.Range("D11").Value = 0.7283763178587769
$ws.Range("E11").Value = 0.853449657483543
$ws.Range("F11").Value = 0.8665861862271473
$ws.Range("G11").Value = 33
